$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1093.0385
$ws.Range("I15").Value = 1093.0385
$ws.Range("K15").Value = 3279.1155
$ws.Range("M15").Value = -3110.1155

$ws.Range("H18").Value = 693.5
$ws.Range("I18").Value = 693.5
$ws.Range("K18").Value = 693.5
$ws.Range("M18").Value = -409.5

$ws.Range("H41").Value = 1626.0358
$ws.Range("I41").Value = 2127.0588
$ws.Range("J41").Value = 851.7273
$ws.Range("K41").Value = 2127.0588
$ws.Range("L41").Value = 851.7273
$ws.Range("M41").Value = -1687.0588
$ws.Range("N41").Value = -1731.7273

$ws.Range("H92").Value = 10005958
$ws.Range("I92").Value = 14293281
$ws.Range("K92").Value = 14293281
$ws.Range("M92").Value = -14292033

$ws.Range("H100").Value = 2322.75
$ws.Range("I100").Value = 2226
$ws.Range("J100").Value = 3000
$ws.Range("K100").Value = 2226
$ws.Range("L100").Value = 3000
$ws.Range("M100").Value = -1685
$ws.Range("N100").Value = -4082

$ws.Range("H135").Value = 17550182
$ws.Range("I135").Value = 18524808
$ws.Range("K135").Value = 166723272
$ws.Range("M135").Value = -166720737

$ws.Range("H138").Value = 4361.271
$ws.Range("I138").Value = 1603.0834
$ws.Range("K138").Value = 4809.2502
$ws.Range("M138").Value = 330.7497999999996

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H34").Value = 30000
$ws.Range("J34").Value = 30000
$ws.Range("L34").Value = 30000
$ws.Range("N34").Value = -30542

$ws.Range("H74").Value = 18525636
$ws.Range("I74").Value = 3804.7273
$ws.Range("J74").Value = 100021700
$ws.Range("K74").Value = 3804.7273
$ws.Range("L74").Value = 100021700
$ws.Range("M74").Value = -2930.7273
$ws.Range("N74").Value = -100023448

$ws.Range("H77").Value = 18525636
$ws.Range("I77").Value = 3804.7273
$ws.Range("J77").Value = 100021700
$ws.Range("K77").Value = 19023.6365
$ws.Range("L77").Value = 500108500
$ws.Range("M77").Value = -14655.6365
$ws.Range("N77").Value = -500117236

$ws.Range("H132").Value = 1711.2909
$ws.Range("I132").Value = 1226.7916
$ws.Range("K132").Value = 3680.3748
$ws.Range("M132").Value = -1150.3748

$ws.Range("H135").Value = 45277.5
$ws.Range("J135").Value = 45277.5
$ws.Range("L135").Value = 45277.5
$ws.Range("N135").Value = -55417.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1045.1724
$ws.Range("I94").Value = 1014.5909
$ws.Range("J94").Value = 1141.2858
$ws.Range("K94").Value = 1014.5909
$ws.Range("L94").Value = 1141.2858
$ws.Range("M94").Value = -563.5909
$ws.Range("N94").Value = -2043.2858

$ws.Range("H105").Value = 5450.4287
$ws.Range("I105").Value = 2329.5
$ws.Range("K105").Value = 2329.5
$ws.Range("M105").Value = -582.5

$ws.Range("H134").Value = 1542.4736
$ws.Range("I134").Value = 1288.9697
$ws.Range("K134").Value = 3866.9091
$ws.Range("M134").Value = -1331.9091

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 624679.8
$ws.Range("I4").Value = 30849.75
$ws.Range("K4").Value = 30849.75
$ws.Range("M4").Value = -30737.75

$ws.Range("H134").Value = 11810.538
$ws.Range("I134").Value = 12227.571
$ws.Range("K134").Value = 36682.713
$ws.Range("M134").Value = -34147.713

$ws.Range("H141").Value = 47360
$ws.Range("J141").Value = 47360
$ws.Range("L141").Value = 47360
$ws.Range("N141").Value = -57720

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 45505344
$ws.Range("I4").Value = 52690084
$ws.Range("K4").Value = 158070252
$ws.Range("M4").Value = -158070140

$ws.Range("H41").Value = 730
$ws.Range("I41").Value = 95
$ws.Range("K41").Value = 285
$ws.Range("M41").Value = 53

$ws.Range("H109").Value = 239222
$ws.Range("I109").Value = 317629.34
$ws.Range("K109").Value = 952888.02
$ws.Range("M109").Value = -951848.02

$ws.Range("H140").Value = 3314.4119
$ws.Range("J140").Value = 4300
$ws.Range("L140").Value = 12900
$ws.Range("N140").Value = -23260

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 2450
$ws.Range("I5").Value = 2450
$ws.Range("K5").Value = 2450
$ws.Range("M5").Value = -2338

$ws.Range("H18").Value = 37076704
$ws.Range("J18").Value = 100000
$ws.Range("L18").Value = 100000
$ws.Range("N18").Value = -100586

$ws.Range("H80").Value = 2303.2727
$ws.Range("I80").Value = 2320.8572
$ws.Range("J80").Value = 2272.5
$ws.Range("K80").Value = 2320.8572
$ws.Range("L80").Value = 2272.5
$ws.Range("M80").Value = -1322.8572
$ws.Range("N80").Value = -4268.5

$ws.Range("H83").Value = 2303.2727
$ws.Range("I83").Value = 2320.8572
$ws.Range("J83").Value = 2272.5
$ws.Range("K83").Value = 11604.286
$ws.Range("L83").Value = 11362.5
$ws.Range("M83").Value = -6612.286
$ws.Range("N83").Value = -21346.5

$ws.Range("H104").Value = 105000
$ws.Range("J104").Value = 105000
$ws.Range("L104").Value = 105000
$ws.Range("N104").Value = -111988

$ws.Range("H123").Value = 25281.125
$ws.Range("J123").Value = 25058.166
$ws.Range("L123").Value = 25058.166
$ws.Range("N123").Value = -29958.166

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 5200000
$ws.Range("J2").Value = 6900000
$ws.Range("L2").Value = 6900000
$ws.Range("N2").Value = -6900224

$ws.Range("H22").Value = 1963
$ws.Range("I22").Value = 1521
$ws.Range("K22").Value = 1521
$ws.Range("M22").Value = -1226

$ws.Range("H27").Value = 1963
$ws.Range("I27").Value = 1521
$ws.Range("K27").Value = 1521
$ws.Range("M27").Value = -1414

$ws.Range("H55").Value = 1367.2
$ws.Range("I55").Value = 1634.25
$ws.Range("K55").Value = 1634.25
$ws.Range("M55").Value = -1461.25

$ws.Range("N74").ClearContents()
$ws.Range("H74").Value = 30000
$ws.Range("I74").Value = 30000
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 30000
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -29002

$ws.Range("N77").ClearContents()
$ws.Range("H77").Value = 30000
$ws.Range("I77").Value = 30000
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 90000
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -85008

$ws.Range("H127").Value = 153142.8
$ws.Range("J127").Value = 153142.8
$ws.Range("L127").Value = 153142.8
$ws.Range("N127").Value = -163062.8

$ws.Range("H132").Value = 3391.8823
$ws.Range("I132").Value = 2656.561
$ws.Range("K132").Value = 7969.683000000001
$ws.Range("M132").Value = -5439.683000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1839.3
$ws.Range("I113").Value = 1849.125
$ws.Range("K113").Value = 5547.375
$ws.Range("M113").Value = -3377.375

$ws.Range("H126").Value = 2685.2415
$ws.Range("I126").Value = 2264.0435
$ws.Range("J126").Value = 4299.8335
$ws.Range("K126").Value = 6792.130500000001
$ws.Range("L126").Value = 12899.5005
$ws.Range("M126").Value = -4322.130500000001
$ws.Range("N126").Value = -17839.5005

$ws.Range("H138").Value = 88614
$ws.Range("J138").Value = 88614
$ws.Range("L138").Value = 88614
$ws.Range("N138").Value = -98894
